# Updated cryptos list on Wed Jan 17 03:56:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.863.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'2.572.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'314.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'99.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.41%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'36.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'7.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "'15.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").Value = "'2.564.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "'42.965.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'6.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "'0.0₃0967"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'69.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'249.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'27.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("D29").Value = "'40.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'10.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "'158.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "'5.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "'3.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "'2.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "'0.0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'18.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  +9.24%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "'23.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +8.35%  "
$ws.Range("D43").Value = "'0.0303"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'3.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "'2.007.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'2.821.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").Value = "'74.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "'81.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.86%  "
